$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered - same style as H1) onto the
# two new header cells, then set their text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF) for rows 2-16.
$data = @{
    2  = @(6, 7)
    3  = @(8, 9)
    4  = @(7, 8)
    5  = @(6, 6)
    6  = @(5, 6)
    7  = @(6, 6)
    8  = @(7, 7)
    9  = @(9, 9)
    10 = @(8, 8)
    11 = @(9, 9)
    12 = @(4, 4)
    13 = @(7, 7)
    14 = @(8, 8)
    15 = @(9, 9)
    16 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
